$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at position 271, pushing existing rows 271:343 down to 273:345
$ws.Rows("271:272").Insert()

# Row 271 (new record)
$ws.Range("A271").Value = 3
$ws.Range("B271").Value = "Femacal de La Calera"
$ws.Range("C271").Value = "Coquimbo"
$ws.Range("D271").Value = 44754
$ws.Range("D271").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E271").Value = 5
$ws.Range("F271").Value = 100112013
$ws.Range("G271").Value = "Alcachofa"
$ws.Range("H271").Value = "Argentina(o)"
$ws.Range("I271").Value = "Primera"
$ws.Range("J271").Value = 105
$ws.Range("K271").Value = 14000
$ws.Range("L271").Value = 15000
$ws.Range("M271").Value = 14524
$ws.Range("N271").Value = "`$/caja 50 unidades"
$ws.Range("O271").Value = "Provincia de Limar" + [char]0x00ED
$ws.Range("P271").Value = 290
$ws.Range("Q271").Value = 50
$ws.Range("R271").Value = "Hortaliza"

# Row 272 (new record)
$ws.Range("A272").Value = 3
$ws.Range("B272").Value = "Femacal de La Calera"
$ws.Range("C272").Value = "Coquimbo"
$ws.Range("D272").Value = 44754
$ws.Range("D272").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E272").Value = 5
$ws.Range("F272").Value = 100112013
$ws.Range("G272").Value = "Alcachofa"
$ws.Range("H272").Value = "Espa" + [char]0x00F1 + "ola"
$ws.Range("I272").Value = "Extra"
$ws.Range("J272").Value = 85
$ws.Range("K272").Value = 17000
$ws.Range("L272").Value = 18000
$ws.Range("M272").Value = 17529
$ws.Range("N272").Value = "`$/caja 30 unidades"
$ws.Range("O272").Value = "Provincia de Limar" + [char]0x00ED
$ws.Range("P272").Value = 584
$ws.Range("Q272").Value = 30
$ws.Range("R272").Value = "Hortaliza"
